$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1441.3636
$ws.Range("I28").Value = 1514
$ws.Range("K28").Value = 1514
$ws.Range("M28").Value = -1029
$ws.Range("H112").Value = 23811218
$ws.Range("I112").Value = 142857870
$ws.Range("J112").Value = 1888.8857
$ws.Range("K112").Value = 428573610
$ws.Range("L112").Value = 5666.6571
$ws.Range("M112").Value = -428572502
$ws.Range("N112").Value = -7882.6571
$ws.Range("H132").Value = 846660.75
$ws.Range("I132").Value = 1654.7451
$ws.Range("J132").Value = 7003133
$ws.Range("K132").Value = 4964.2353
$ws.Range("L132").Value = 21009399
$ws.Range("M132").Value = -2434.2353
$ws.Range("N132").Value = -21014459
$ws.Range("H137").Value = 1667956.1
$ws.Range("I137").Value = 2326581.8
$ws.Range("K137").Value = 6979745.399999999
$ws.Range("M137").Value = -6977195.399999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 22267770
$ws.Range("I61").Value = 25026064
$ws.Range("J61").Value = 201422.8
$ws.Range("K61").Value = 25026064
$ws.Range("L61").Value = 201422.8
$ws.Range("M61").Value = -25025852
$ws.Range("N61").Value = -201846.8
$ws.Range("H136").Value = 22267770
$ws.Range("I136").Value = 25026064
$ws.Range("J136").Value = 201422.8
$ws.Range("K136").Value = 75078192
$ws.Range("L136").Value = 604268.3999999999
$ws.Range("M136").Value = -75075642
$ws.Range("N136").Value = -609368.3999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2172.0312
$ws.Range("I134").Value = 1274.1666
$ws.Range("K134").Value = 3822.4998
$ws.Range("M134").Value = -1287.4998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13250
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 13250
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 13250
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -13840
$ws.Range("H34").Value = 13250
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 13250
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 13250
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -13654
$ws.Range("H52").Value = 38420
$ws.Range("J52").Value = 38420
$ws.Range("L52").Value = 38420
$ws.Range("N52").Value = -39008
$ws.Range("H58").Value = 18183158
$ws.Range("I58").Value = 20409376
$ws.Range("J58").Value = 2378.5
$ws.Range("K58").Value = 20409376
$ws.Range("L58").Value = 2378.5
$ws.Range("M58").Value = -20409173
$ws.Range("N58").Value = -2784.5
$ws.Range("H132").Value = 39171.113
$ws.Range("I132").Value = 23604.2
$ws.Range("J132").Value = 126735
$ws.Range("K132").Value = 70812.60000000001
$ws.Range("L132").Value = 380205
$ws.Range("M132").Value = -68282.60000000001
$ws.Range("N132").Value = -385265
$ws.Range("H134").Value = 32313.97
$ws.Range("I134").Value = 1429.8077
$ws.Range("J134").Value = 121534.89
$ws.Range("K134").Value = 4289.4231
$ws.Range("L134").Value = 364604.67
$ws.Range("M134").Value = -1754.4231
$ws.Range("N134").Value = -369674.67
$ws.Range("H136").Value = 18183158
$ws.Range("I136").Value = 20409376
$ws.Range("J136").Value = 2378.5
$ws.Range("K136").Value = 61228128
$ws.Range("L136").Value = 7135.5
$ws.Range("M136").Value = -61225578
$ws.Range("N136").Value = -12235.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 440
$ws.Range("I17").Value = 550
$ws.Range("J17").Value = 366.66666
$ws.Range("K17").Value = 1650
$ws.Range("L17").Value = 1099.99998
$ws.Range("M17").Value = -1481
$ws.Range("N17").Value = -1437.99998
$ws.Range("H97").Value = 2615.889
$ws.Range("I97").Value = 4395.8
$ws.Range("J97").Value = 391
$ws.Range("K97").Value = 13187.4
$ws.Range("L97").Value = 1173
$ws.Range("M97").Value = -12691.4
$ws.Range("N97").Value = -2165
$ws.Range("H129").Value = 3207155.5
$ws.Range("I129").Value = 2073.3333
$ws.Range("J129").Value = 4168680
$ws.Range("K129").Value = 6219.999899999999
$ws.Range("L129").Value = 12506040
$ws.Range("M129").Value = -1219.999899999999
$ws.Range("N129").Value = -12516040
$ws.Range("H130").Value = 2745
$ws.Range("J130").Value = 3181.25
$ws.Range("L130").Value = 9543.75
$ws.Range("N130").Value = -19583.75
$ws.Range("H131").Value = 1102.4286
$ws.Range("J131").Value = 1286.6562
$ws.Range("L131").Value = 3859.9686
$ws.Range("N131").Value = -13939.9686
$ws.Range("H136").Value = 2665.8948
$ws.Range("I136").Value = 1881.125
$ws.Range("J136").Value = 3236.6365
$ws.Range("K136").Value = 5643.375
$ws.Range("L136").Value = 9709.9095
$ws.Range("M136").Value = -543.375
$ws.Range("N136").Value = -19909.9095
$ws.Range("H140").Value = 2559.6562
$ws.Range("I140").Value = 3039.0908
$ws.Range("J140").Value = 2308.524
$ws.Range("K140").Value = 9117.2724
$ws.Range("L140").Value = 6925.572
$ws.Range("M140").Value = -3937.2724
$ws.Range("N140").Value = -17285.572

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 11666774
$ws.Range("I14").Value = 11666774
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 11666774
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -11666606
$ws.Range("N14").ClearContents()
$ws.Range("H107").Value = 324.70587
$ws.Range("J107").Value = 572.5
$ws.Range("L107").Value = 572.5
$ws.Range("N107").Value = -4412.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H22").Value = 758.2258
$ws.Range("I22").Value = 496.27274
$ws.Range("J22").Value = 902.3
$ws.Range("K22").Value = 496.27274
$ws.Range("L22").Value = 902.3
$ws.Range("M22").Value = -201.27274
$ws.Range("N22").Value = -1492.3
$ws.Range("H27").Value = 758.2258
$ws.Range("I27").Value = 496.27274
$ws.Range("J27").Value = 902.3
$ws.Range("K27").Value = 496.27274
$ws.Range("L27").Value = 902.3
$ws.Range("M27").Value = -389.27274
$ws.Range("N27").Value = -1116.3
$ws.Range("H46").Value = 726.8889
$ws.Range("I46").Value = 656.35297
$ws.Range("J46").Value = 790
$ws.Range("K46").Value = 656.35297
$ws.Range("L46").Value = 790
$ws.Range("M46").Value = -468.35297
$ws.Range("N46").Value = -1166
$ws.Range("H88").Value = 36731.05
$ws.Range("J88").Value = 36731.05
$ws.Range("L88").Value = 36731.05
$ws.Range("N88").Value = -37587.05
$ws.Range("H91").Value = 36731.05
$ws.Range("J91").Value = 36731.05
$ws.Range("L91").Value = 36731.05
$ws.Range("N91").Value = -39695.05
$ws.Range("H132").Value = 83701.8
$ws.Range("I132").Value = 69841.664
$ws.Range("J132").Value = 104492
$ws.Range("K132").Value = 209524.992
$ws.Range("L132").Value = 313476
$ws.Range("M132").Value = -206994.992
$ws.Range("N132").Value = -318536
$ws.Range("H136").Value = 70247
$ws.Range("I136").Value = 61932.668
$ws.Range("J136").Value = 81759.16
$ws.Range("K136").Value = 185798.004
$ws.Range("L136").Value = 245277.48
$ws.Range("M136").Value = -183248.004
$ws.Range("N136").Value = -250377.48

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 486.8889
$ws.Range("I107").Value = 417.5
$ws.Range("K107").Value = 1252.5
$ws.Range("M107").Value = 667.5
$ws.Range("H132").Value = 48372.125
$ws.Range("I132").Value = 37016.508
$ws.Range("J132").Value = 110828
$ws.Range("K132").Value = 111049.524
$ws.Range("L132").Value = 332484
$ws.Range("M132").Value = -108519.524
$ws.Range("N132").Value = -337544
$ws.Range("H136").Value = 30263.984
$ws.Range("I136").Value = 17959.14
$ws.Range("K136").Value = 53877.42
$ws.Range("M136").Value = -51327.42
